# Generate Report for Handback
# Update the Correspond Handoff / Handback timestamps for the
# 15f878de-... entry on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 11:37:32"
$wsZhCn.Range("H3").Value = "2016-03-24 11:38:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 11:37:37"
$wsDeDe.Range("H3").Value = "2016-03-24 11:38:15"
